$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "number of worked hours" for the 28.10.2022 and 04.11.2022
#     rows (rows 10-13) from 6 to 7. ---
$ws.Range("D10").Value = 7
$ws.Range("D11").Value = 7
$ws.Range("D12").Value = 7
$ws.Range("D13").Value = 7

# --- Append the new journal entries (rows 14-18). ---
# "11.11.2022" is ambiguous (day=11, month=11), so a plain assignment gets
# auto-converted into a date serial by the smart-entry heuristic. Route it
# through a scratch formula cell + copy/paste-values instead, which keeps it
# as literal text without leaving any quote-prefix / number-format residue
# behind in the cell's style.
$ws.Range("ZZ1").Formula = "=""11.11.2022"""
$ws.Range("ZZ1").Copy()
$ws.Range("A14").PasteSpecial(-4163)
$ws.Range("A15").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

$ws.Range("B14").Value = "Sandro"
$ws.Range("C15").Value = "Discussed what final project we could do and looked online for the materials."
$ws.Range("C14").Value = "Discussed what final project we could do and looked online for the materials"
$ws.Range("D14").Value = 5

$ws.Range("B15").Value = "Alexandru "
$ws.Range("D15").Value = 5

$ws.Range("A16").Value = "15.11.2022"
$ws.Range("B16").Value = "Alexandru "
$ws.Range("C16").Value = "Worked on the tutorial"
$ws.Range("D16").Value = 1

$ws.Range("A17").Value = "16.11.2023"
$ws.Range("B17").Value = "Alexandru "
$ws.Range("C17").Value = "Worked on the tutorial"
$ws.Range("D17").Value = 3

$ws.Range("A18").Value = "16.11.2024"
$ws.Range("B18").Value = "Sandro "
$ws.Range("C18").Value = "Worked on the tutorial"
$ws.Range("D18").Value = "??"

$ws.Range("D21").Select()
